# Updates odds values on rows 2, 4 and 8 of Sheet1 as per the
# "Atualizando o arquivo XLSX" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$row2 = @{
    "H2" = 2.7
    "I2" = 2.5
    "J2" = 4.5
    "K2" = 1.73
    "L2" = 3.6
    "M2" = 1.18
    "N2" = 4.5
    "O2" = 1.83
    "P2" = 1.83
    "S2" = 3.6
    "T2" = 1.29
    "U2" = 6.4
    "V2" = 1.11
    "W2" = 9
    "X2" = 1.07
    "Y2" = 1.83
    "Z2" = 1.98
    "AA2" = 2.75
    "AB2" = 1.4
    "AC2" = 6
    "AH2" = 67
    "AI2" = 4.33
    "AK2" = 26
    "AL2" = 126
    "AM2" = 5
    "AO2" = 12
    "AP2" = 26
    "AQ2" = 34
    "AR2" = 51
}

foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}

# Row 4 updates
$row4 = @{
    "G4" = 7
    "I4" = 1.39
    "J4" = 6.5
    "M4" = 1.02
    "O4" = 1.13
    "U4" = 1.98
    "V4" = 1.88
    "W4" = 2.37
    "X4" = 1.5
    "AA4" = 1.75
    "AB4" = 2
    "AD4" = 41
    "AE4" = 21
    "AF4" = 81
    "AI4" = 15
    "AJ4" = 9
    "AL4" = 51
    "AM4" = 8.5
    "AN4" = 7.5
}

foreach ($addr in $row4.Keys) {
    $ws.Range($addr).Value = $row4[$addr]
}

# Row 8 updates
$row8 = @{
    "M8" = 1.05
    "N8" = 8.5
    "S8" = 1.98
    "T8" = 1.83
    "W8" = 3.4
    "X8" = 1.3
}

foreach ($addr in $row8.Keys) {
    $ws.Range($addr).Value = $row8[$addr]
}
